{"js": "// Update the date line and every \"NNN\u00d7N=\" problem cell in the table to the\n// new values from the latest generated worksheet.\nconst replacements = [\n  [\"2025-12-14 Sunday\", \"2025-12-15 Monday\"],\n  [\"595\u00d73=\", \"137\u00d72=\"],\n  [\"954\u00d75=\", \"538\u00d75=\"],\n  [\"525\u00d77=\", \"155\u00d77=\"],\n  [\"722\u00d75=\", \"906\u00d73=\"],\n  [\"249\u00d74=\", \"205\u00d72=\"],\n  [\"391\u00d75=\", \"226\u00d78=\"],\n  [\"891\u00d79=\", \"111\u00d77=\"],\n  [\"913\u00d77=\", \"152\u00d74=\"],\n  [\"250\u00d76=\", \"718\u00d79=\"],\n  [\"362\u00d78=\", \"561\u00d75=\"],\n  [\"731\u00d77=\", \"166\u00d78=\"],\n  [\"631\u00d77=\", \"882\u00d76=\"],\n  [\"142\u00d74=\", \"306\u00d75=\"],\n  [\"980\u00d79=\", \"556\u00d77=\"],\n  [\"331\u00d73=\", \"476\u00d77=\"],\n  [\"862\u00d78=\", \"365\u00d77=\"],\n  [\"474\u00d76=\", \"964\u00d75=\"],\n  [\"122\u00d74=\", \"559\u00d74=\"],\n  [\"862\u00d79=\", \"541\u00d73=\"],\n  [\"620\u00d72=\", \"279\u00d75=\"],\n  [\"867\u00d76=\", \"994\u00d76=\"],\n  [\"269\u00d78=\", \"640\u00d79=\"],\n  [\"562\u00d73=\", \"550\u00d79=\"],\n  [\"672\u00d74=\", \"555\u00d72=\"],\n  [\"981\u00d79=\", \"407\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"NNN\u00d7N=\" problem cell in the table to the\n# new values from the latest generated worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-12-14 Sunday\", \"2025-12-15 Monday\"),\n  @(\"595\u00d73=\", \"137\u00d72=\"),\n  @(\"954\u00d75=\", \"538\u00d75=\"),\n  @(\"525\u00d77=\", \"155\u00d77=\"),\n  @(\"722\u00d75=\", \"906\u00d73=\"),\n  @(\"249\u00d74=\", \"205\u00d72=\"),\n  @(\"391\u00d75=\", \"226\u00d78=\"),\n  @(\"891\u00d79=\", \"111\u00d77=\"),\n  @(\"913\u00d77=\", \"152\u00d74=\"),\n  @(\"250\u00d76=\", \"718\u00d79=\"),\n  @(\"362\u00d78=\", \"561\u00d75=\"),\n  @(\"731\u00d77=\", \"166\u00d78=\"),\n  @(\"631\u00d77=\", \"882\u00d76=\"),\n  @(\"142\u00d74=\", \"306\u00d75=\"),\n  @(\"980\u00d79=\", \"556\u00d77=\"),\n  @(\"331\u00d73=\", \"476\u00d77=\"),\n  @(\"862\u00d78=\", \"365\u00d77=\"),\n  @(\"474\u00d76=\", \"964\u00d75=\"),\n  @(\"122\u00d74=\", \"559\u00d74=\"),\n  @(\"862\u00d79=\", \"541\u00d73=\"),\n  @(\"620\u00d72=\", \"279\u00d75=\"),\n  @(\"867\u00d76=\", \"994\u00d76=\"),\n  @(\"269\u00d78=\", \"640\u00d79=\"),\n  @(\"562\u00d73=\", \"550\u00d79=\"),\n  @(\"672\u00d74=\", \"555\u00d72=\"),\n  @(\"981\u00d79=\", \"407\u00d79=\")\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $p[0]\n  $find.Replacement.Text = $p[1]\n  $find.Execute($p[0], $false, $false, $false, $false, $false, $true, 1, $false, $p[1], 2) | Out-Null\n}\n"}
